$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.542 (0.013)"
$ws.Range("D2").Value = "0.486 (0.009)"
$ws.Range("E2").Value = "0.217 (0.009)"

$ws.Range("C3").Value = "0.652 (0.010)"
$ws.Range("D3").Value = "0.576 (0.009)"
$ws.Range("E3").Value = "0.226 (0.002)"

$ws.Range("C4").Value = "0.743 (0.013)"
$ws.Range("D4").Value = "0.656 (0.014)"
$ws.Range("E4").Value = "0.237 (0.009)"

$ws.Range("C5").Value = "0.786 (0.007)"
$ws.Range("D5").Value = "0.683 (0.007)"
$ws.Range("E5").Value = "0.242 (0.002)"

$ws.Range("C6").Value = "0.520 (0.016)"
$ws.Range("D6").Value = "0.501 (0.012)"
$ws.Range("E6").Value = "0.212 (0.022)"

$ws.Range("C7").Value = "0.616 (0.008)"
$ws.Range("D7").Value = "0.587 (0.007)"
$ws.Range("E7").Value = "0.231 (0.001)"

$ws.Range("C8").Value = "0.666 (0.019)"
$ws.Range("D8").Value = "0.654 (0.015)"
$ws.Range("E8").Value = "0.243 (0.008)"

$ws.Range("C9").Value = "0.742 (0.009)"
$ws.Range("D9").Value = "0.709 (0.008)"
$ws.Range("E9").Value = "0.248 (0.002)"

$ws.Range("C10").Value = "0.406 (0.018)"
$ws.Range("D10").Value = "0.418 (0.016)"
$ws.Range("E10").Value = "0.113 (0.011)"

$ws.Range("C11").Value = "0.614 (0.010)"
$ws.Range("D11").Value = "0.589 (0.006)"
$ws.Range("E11").Value = "0.306 (0.006)"

$ws.Range("C12").Value = "0.643 (0.026)"
$ws.Range("D12").Value = "0.637 (0.020)"
$ws.Range("E12").Value = "0.276 (0.027)"

$ws.Range("C13").Value = "0.744 (0.006)"
$ws.Range("D13").Value = "0.707 (0.007)"
$ws.Range("E13").Value = "0.406 (0.007)"
